# "catching up on weeks 5,6,7"
#
# - Images sheet: mark Jason, Nate, and Shannon as Eliminated ("Yes") since
#   they were voted out in weeks 5, 6, and 7 respectively.
# - Elimination_Table sheet: append the new elimination rows for
#   Week 5 -> Jason, Week 6 -> Shannon, Week 7 -> Nate.
# - Update the saved cursor/selection on each sheet and make "Images" the
#   active (front-most) tab again.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Images")
$ws2 = $wb.Worksheets.Item("Elimination_Table")

# Mark newly eliminated players on the Images sheet.
$ws1.Range("C5").Value  = "Yes"   # Jason  - Week 5
$ws1.Range("C11").Value = "Yes"   # Nate   - Week 7
$ws1.Range("C16").Value = "Yes"   # Shannon- Week 6

# Record weeks 5-7 on the Elimination_Table sheet.
$ws2.Range("A7").Value = "Week 5"
$ws2.Range("B7").Value = "Jason"
$ws2.Range("A8").Value = "Week 6"
$ws2.Range("B8").Value = "Shannon"
$ws2.Range("A9").Value = "Week 7"
$ws2.Range("B9").Value = "Nate"

# Restore each sheet's own last-used selection...
$ws2.Activate()
$ws2.Range("H18").Select()

# ...then leave "Images" as the active/front tab with its own selection.
$ws1.Activate()
$ws1.Range("C12").Select()
